$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 67.8125
$ws.Range("I9").Value = 65.583336
$ws.Range("K9").Value = 65.583336
$ws.Range("M9").Value = 103.416664

$ws.Range("H17").Value = 4350657
$ws.Range("J17").Value = 4548328
$ws.Range("L17").Value = 13644984
$ws.Range("N17").Value = -13645320

$ws.Range("H33").Value = 638.6111
$ws.Range("I33").Value = 605.82355
$ws.Range("K33").Value = 605.82355
$ws.Range("M33").Value = -376.82355

$ws.Range("H80").Value = 9066.666999999999
$ws.Range("I80").Value = 7464.9287
$ws.Range("K80").Value = 22394.7861
$ws.Range("M80").Value = -21396.7861

$ws.Range("H83").Value = 9066.666999999999
$ws.Range("I83").Value = 7464.9287
$ws.Range("K83").Value = 67184.35830000001
$ws.Range("M83").Value = -62192.35830000001

$ws.Range("H96").Value = 2076.375
$ws.Range("I96").Value = 1456.8334
$ws.Range("J96").Value = 3935
$ws.Range("K96").Value = 4370.5002
$ws.Range("L96").Value = 11805
$ws.Range("M96").Value = -2997.5002
$ws.Range("N96").Value = -14551

$ws.Range("H105").Value = 54010
$ws.Range("J105").Value = 45680.332
$ws.Range("L105").Value = 45680.332
$ws.Range("N105").Value = -52668.332

$ws.Range("H116").Value = 4508.75
$ws.Range("I116").Value = 4129
$ws.Range("J116").Value = 5040.4
$ws.Range("K116").Value = 4129
$ws.Range("L116").Value = 5040.4
$ws.Range("M116").Value = -687
$ws.Range("N116").Value = -11924.4

$ws.Range("H132").Value = 966.8043
$ws.Range("I132").Value = 966.8043
$ws.Range("K132").Value = 2900.4129
$ws.Range("M132").Value = -370.4129000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2137.111
$ws.Range("I2").Value = 1056.3448
$ws.Range("K2").Value = 1056.3448
$ws.Range("M2").Value = -943.3448000000001

$ws.Range("H32").Value = 5651.5654
$ws.Range("I32").Value = 5554.9556
$ws.Range("K32").Value = 5554.9556
$ws.Range("M32").Value = -5267.9556

$ws.Range("H102").Value = 4990.4165
$ws.Range("I102").Value = 3765
$ws.Range("J102").Value = 8666.666999999999
$ws.Range("K102").Value = 3765
$ws.Range("L102").Value = 8666.666999999999
$ws.Range("M102").Value = -2143
$ws.Range("N102").Value = -11910.667

$ws.Range("H116").Value = 2137.111
$ws.Range("I116").Value = 1056.3448
$ws.Range("K116").Value = 1056.3448
$ws.Range("M116").Value = 1237.6552

$ws.Range("H132").Value = 2749.5908
$ws.Range("I132").Value = 2665.0667
$ws.Range("K132").Value = 7995.2001
$ws.Range("M132").Value = -5465.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2137.111
$ws.Range("I3").Value = 1056.3448
$ws.Range("K3").Value = 1056.3448
$ws.Range("M3").Value = -942.3448000000001

$ws.Range("H94").Value = 1072.4828
$ws.Range("I94").Value = 1072.4828
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1072.4828
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -621.4828
$ws.Range("N94").ClearContents()

$ws.Range("H108").Value = 58999.668
$ws.Range("I108").Value = 58999
$ws.Range("K108").Value = 58999
$ws.Range("M108").Value = -55159

$ws.Range("H134").Value = 2351.5
$ws.Range("I134").Value = 2266.913
$ws.Range("K134").Value = 6800.739
$ws.Range("M134").Value = -4265.739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30307272
$ws.Range("I31").Value = 43481096
$ws.Range("J31").Value = 7474.8
$ws.Range("K31").Value = 43481096
$ws.Range("L31").Value = 7474.8
$ws.Range("M31").Value = -43480801
$ws.Range("N31").Value = -8064.8

$ws.Range("H34").Value = 30307272
$ws.Range("I34").Value = 43481096
$ws.Range("J34").Value = 7474.8
$ws.Range("K34").Value = 43481096
$ws.Range("L34").Value = 7474.8
$ws.Range("M34").Value = -43480894
$ws.Range("N34").Value = -7878.8

$ws.Range("H86").Value = 8182.5713
$ws.Range("J86").Value = 6093.6
$ws.Range("L86").Value = 6093.6
$ws.Range("N86").Value = -8339.6

$ws.Range("H89").Value = 8182.5713
$ws.Range("J89").Value = 6093.6
$ws.Range("L89").Value = 30468
$ws.Range("N89").Value = -41700

$ws.Range("H99").Value = 5793.7036
$ws.Range("I99").Value = 6198.9287
$ws.Range("J99").Value = 5357.3076
$ws.Range("K99").Value = 6198.9287
$ws.Range("L99").Value = 5357.3076
$ws.Range("M99").Value = -4700.9287
$ws.Range("N99").Value = -8353.3076

$ws.Range("H126").Value = 5793.7036
$ws.Range("I126").Value = 6198.9287
$ws.Range("J126").Value = 5357.3076
$ws.Range("K126").Value = 18596.7861
$ws.Range("L126").Value = 16071.9228
$ws.Range("M126").Value = -16126.7861
$ws.Range("N126").Value = -21011.9228

$ws.Range("H141").Value = 213355.28
$ws.Range("J141").Value = 235881.25
$ws.Range("L141").Value = 235881.25
$ws.Range("N141").Value = -246241.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 35721536
$ws.Range("J68").Value = 8375.75
$ws.Range("L68").Value = 25127.25
$ws.Range("N68").Value = -26749.25

$ws.Range("H71").Value = 35721536
$ws.Range("J71").Value = 8375.75
$ws.Range("L71").Value = 75381.75
$ws.Range("N71").Value = -83493.75

$ws.Range("H121").Value = 67573.53
$ws.Range("J121").Value = 111987.89
$ws.Range("L121").Value = 335963.67
$ws.Range("N121").Value = -338583.67

$ws.Range("H131").Value = 16669001
$ws.Range("I131").Value = 31251076
$ws.Range("J131").Value = 3772.7856
$ws.Range("K131").Value = 93753228
$ws.Range("L131").Value = 11318.3568
$ws.Range("M131").Value = -93748188
$ws.Range("N131").Value = -21398.3568

$ws.Range("H137").Value = 2421.3333
$ws.Range("I137").Value = 2121.8462
$ws.Range("K137").Value = 6365.5386
$ws.Range("M137").Value = -1265.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3044.52
$ws.Range("J80").Value = 3805.6667
$ws.Range("L80").Value = 3805.6667
$ws.Range("N80").Value = -5801.6667

$ws.Range("H83").Value = 3044.52
$ws.Range("J83").Value = 3805.6667
$ws.Range("L83").Value = 19028.3335
$ws.Range("N83").Value = -29012.3335

$ws.Range("H126").Value = 4960.0415
$ws.Range("I126").Value = 4545.0835
$ws.Range("J126").Value = 5375
$ws.Range("K126").Value = 13635.2505
$ws.Range("L126").Value = 16125
$ws.Range("M126").Value = -11165.2505
$ws.Range("N126").Value = -21065

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1370
$ws.Range("I22").Value = 1890
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 1890
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = -1595
$ws.Range("N22").Value = -1440

$ws.Range("H27").Value = 1370
$ws.Range("I27").Value = 1890
$ws.Range("J27").Value = 850
$ws.Range("K27").Value = 1890
$ws.Range("L27").Value = 850
$ws.Range("M27").Value = -1783
$ws.Range("N27").Value = -1064

$ws.Range("H40").Value = 3867.611
$ws.Range("J40").Value = 3223.25
$ws.Range("L40").Value = 3223.25
$ws.Range("N40").Value = -3495.25

$ws.Range("H46").Value = 6982.1064
$ws.Range("I46").Value = 6554
$ws.Range("K46").Value = 6554
$ws.Range("M46").Value = -6366

$ws.Range("H68").Value = 82378.46000000001
$ws.Range("J68").Value = 131497.25
$ws.Range("L68").Value = 131497.25
$ws.Range("N68").Value = -132995.25

$ws.Range("H71").Value = 82378.46000000001
$ws.Range("J71").Value = 131497.25
$ws.Range("L71").Value = 657486.25
$ws.Range("N71").Value = -664974.25

$ws.Range("H132").Value = 208780.77
$ws.Range("I132").Value = 230610.8
$ws.Range("K132").Value = 691832.3999999999
$ws.Range("M132").Value = -689302.3999999999

$ws.Range("H136").Value = 5103.486
$ws.Range("I136").Value = 5370.7334
$ws.Range("K136").Value = 16112.2002
$ws.Range("M136").Value = -13562.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H126").Value = 3643
$ws.Range("J126").Value = 5785.4287
$ws.Range("L126").Value = 17356.2861
$ws.Range("N126").Value = -22296.2861

$ws.Range("H132").Value = 149634.78
$ws.Range("I132").Value = 166294.7
$ws.Range("K132").Value = 498884.1
$ws.Range("M132").Value = -496354.1
